$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Range("H6").Value2 = 17637
$ws.Range("S6").Value2 = 402211

# Row 8
$ws.Range("H8").Value2 = 20828
$ws.Range("I8").Value2 = 4509
$ws.Range("S8").Value2 = 563375

# Row 17
$ws.Range("O17").Value2 = 266
$ws.Range("P17").Value2 = 88
$ws.Range("Q17").Value2 = 171347

# Row 18
$ws.Range("I18").Value2 = 5950
$ws.Range("O18").Value2 = 118
$ws.Range("P18").Value2 = 88
$ws.Range("Q18").Value2 = 54453

# Row 30
$ws.Range("P30").Value2 = 80

# Row 31
$ws.Range("H31").Value2 = 7508
$ws.Range("O31").Value2 = 10
$ws.Range("Q31").Value2 = 19320
$ws.Range("S31").Value2 = 494225

# Row 32
$ws.Range("D32").Value2 = 5340
$ws.Range("G32").Value2 = 11377
$ws.Range("H32").Value2 = 24928
$ws.Range("I32").Value2 = 4486
$ws.Range("S32").Value2 = 554641

# Row 35
$ws.Range("H35").Value2 = 10635

# Row 37
$ws.Range("I37").Value2 = 2895

# Row 50
$ws.Range("D50").Value2 = 5368
$ws.Range("F50").Value2 = 9213
$ws.Range("G50").Value2 = 9201
$ws.Range("H50").Value2 = 27745
$ws.Range("O50").Value2 = 72
$ws.Range("Q50").Value2 = 143512
